$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the bank account details next to their labels in column D
$ws.Range("D5").Value = "Kashif Jamal Ahmed"
$ws.Range("D6").Value = "Emirates NBD"
$ws.Range("D7").Value = "Ibn Batutta Mall, Dubai"

# IBAN entered before the account number (matches the author's original fill order)
$ws.Range("D9").Value = "AE96 0260 0002 1580 6638 701"

# Account number keeps a leading zero, so format it as text-like "00000" numbers
$ws.Range("D8").NumberFormat = "00000"
$ws.Range("D8").HorizontalAlignment = -4131
$ws.Range("D8").VerticalAlignment = -4108
$ws.Range("D8").Value = "02 1580 6638 701"

$ws.Range("D10").Value = "EBILAEAD"

# Keep the focus on D10, matching the saved selection in the workbook
$ws.Activate()
$ws.Range("D10").Select()
$excel.ActiveWindow.Zoom = 215
